$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new value would otherwise be auto-converted to a number,
# so they stay text strings like the rest of the Price column (matches original inlineStr formatting).
$textCells = @('D4', 'D5', 'D6', 'D7', 'D8', 'D9', 'D10', 'D11', 'D14', 'D15', 'D18', 'D19', 'D21', 'D23', 'D25', 'D26', 'D27', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D37', 'D38', 'D41', 'D42', 'D43', 'D44', 'D46', 'D48', 'D49', 'D50', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin values row by row.
$ws.Range('D2').Value = '34.628.08'
$ws.Range('E2').Value = '  +0.75%  '

$ws.Range('D3').Value = '1.819.34'
$ws.Range('E3').Value = '  +1.03%  '

$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.19%  '

$ws.Range('D5').Value = '225.83'
$ws.Range('E5').Value = '  +0.31%  '

$ws.Range('D6').Value = '0.607'
$ws.Range('E6').Value = '  +1.30%  '

$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.19%  '

$ws.Range('D8').Value = '44.76'
$ws.Range('E8').Value = '  +23.56%  '

$ws.Range('D9').Value = '0.296'
$ws.Range('E9').Value = '  +1.38%  '

$ws.Range('D10').Value = '0.0679'
$ws.Range('E10').Value = '  +0.33%  '

$ws.Range('D11').Value = '0.0999'
$ws.Range('E11').Value = '  +3.53%  '

$ws.Range('D12').Value = '2.082.17'
$ws.Range('E12').Value = '  +1.03%  '

$ws.Range('D13').Value = '1.823.79'
$ws.Range('E13').Value = '  +0.88%  '

$ws.Range('D14').Value = '11.16'
$ws.Range('E14').Value = '  -0.98%  '

$ws.Range('D15').Value = '0.641'
$ws.Range('E15').Value = '  +2.18%  '

$ws.Range('E16').Value = '  +1.41%  '

$ws.Range('D17').Value = '34.575.75'
$ws.Range('E17').Value = '  +0.58%  '

$ws.Range('D18').Value = '68.17'
$ws.Range('E18').Value = '  -0.41%  '

$ws.Range('D19').Value = '242.22'
$ws.Range('E19').Value = '  +0.00%  '

$ws.Range('D20').Value = '0.0₃0781'
$ws.Range('E20').Value = '  +1.31%  '

$ws.Range('D21').Value = '11.71'
$ws.Range('E21').Value = '  +4.40%  '

$ws.Range('E22').Value = '  -0.03%  '

$ws.Range('D23').Value = '4.50'
$ws.Range('E23').Value = '  +10.36%  '

$ws.Range('E24').Value = '  -2.38%  '

$ws.Range('D25').Value = '170.78'
$ws.Range('E25').Value = '  -0.06%  '

$ws.Range('D26').Value = '7.83'
$ws.Range('E26').Value = '  -0.49%  '

$ws.Range('D27').Value = '17.69'
$ws.Range('E27').Value = '  +2.10%  '

$ws.Range('E28').Value = '  +0.63%  '

$ws.Range('D29').Value = '0.998'
$ws.Range('E29').Value = '  -0.21%  '

$ws.Range('D30').Value = '3.86'
$ws.Range('E30').Value = '  +1.88%  '

$ws.Range('D31').Value = '1.24'
$ws.Range('E31').Value = '  +1.03%  '

$ws.Range('D32').Value = '3.93'
$ws.Range('E32').Value = '  +0.76%  '

$ws.Range('D33').Value = '0.0520'
$ws.Range('E33').Value = '  +1.17%  '

$ws.Range('D34').Value = '1.84'
$ws.Range('E34').Value = '  +3.25%  '

$ws.Range('E35').Value = '  +11.34%  '

$ws.Range('E36').Value = '  +1.26%  '

$ws.Range('B37').Value = 'InjectiveProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D37').Value = '15.40'
$ws.Range('E37').Value = '  +16.05%  '

$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').Value = '2.43'
$ws.Range('E38').Value = '  +3.39%  '

$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').Value = '1.321.96'
$ws.Range('E39').Value = '  -2.83%  '

$ws.Range('E40').Value = '  +0.23%  '

$ws.Range('D41').Value = '0.0191'
$ws.Range('E41').Value = '  +3.12%  '

$ws.Range('B42').Value = 'WEMIXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D42').Value = '1.22'
$ws.Range('E42').Value = '  +5.07%  '

$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').Value = '0.957'
$ws.Range('E43').Value = '  +2.37%  '

$ws.Range('D44').Value = '2.81'
$ws.Range('E44').Value = '  +1.51%  '

$ws.Range('E45').Value = '  -0.16%  '

$ws.Range('D46').Value = '0.0518'
$ws.Range('E46').Value = '  +3.94%  '

$ws.Range('D47').Value = '1.981.35'
$ws.Range('E47').Value = '  +0.96%  '

$ws.Range('D48').Value = '5.89'
$ws.Range('E48').Value = '  +2.04%  '

$ws.Range('D49').Value = '0.999'
$ws.Range('E49').Value = '  -0.14%  '

$ws.Range('D50').Value = '101.55'
$ws.Range('E50').Value = '  -0.70%  '

$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D51').Value = '48.63'
$ws.Range('E51').Value = '  +0.85%  '

